# Change the "MapRendererDlgInvisibleOverlay" / "使用数字渲染透明覆盖物" row
# into a new "MapRendererDlgInvisibleInGame" / "显示游戏内不可见对象" row,
# moved to the end of the table (mirrors the author's commit: renaming the
# map renderer's "numeric overlay" option to "invisible in game").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row that currently holds the old key so this keeps working even
# if the sheet layout shifts a bit.
$finder = $ws.Range("A:A").Find("MapRendererDlgInvisibleOverlay")
if ($finder -ne $null) {
    $oldRow = $finder.Row
} else {
    $oldRow = 792
}

# Remove that entire row - every row below it shifts up by one, exactly like
# a user selecting the row and deleting it.
$ws.Rows($oldRow).Delete()

# The sheet's last used row (now one row smaller) gets the brand new entry
# appended after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "MapRendererDlgInvisibleInGame"
$ws.Cells.Item($newRow, 2).Value = "显示游戏内不可见对象"

# Match the formatting used by the rest of the table (same style as the row
# directly above it) without disturbing the values we just wrote.
$srcFmt = $ws.Range("A" + ($newRow - 1) + ":B" + ($newRow - 1))
$dstFmt = $ws.Range("A" + $newRow + ":B" + $newRow)
$srcFmt.Copy()
$dstFmt.PasteSpecial(-4122)

$ws.Range("A" + $newRow).Select()
